$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the dates: D2/D3 (currently 44559) -> 44574, D4/D5 (currently 44574) -> 44559
$ws.Range("D2").Value = 44574
$ws.Range("D3").Value = 44574
$ws.Range("D4").Value = 44559
$ws.Range("D5").Value = 44559
